$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password value in C2 from "1234Test" to "1234Tests"
$ws.Range("C2").Value = "1234Tests"

# Move the active selection from D12 to E8
$ws.Range("E8").Select()
